# Update odds values for Jogos_da_Semana_FlashScore_2025-03-05.xlsx
# per the FlashScore odds refresh commit "Atualizando o arquivo XLSX".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1.55
$ws.Range("H3").Value = 3.45
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 2.07
$ws.Range("L3").Value = 6.6
$ws.Range("N3").Value = 2.67
$ws.Range("O3").Value = 2.2
$ws.Range("P3").Value = 1.6
$ws.Range("Q3").Value = 3.85
$ws.Range("R3").Value = 1.22
$ws.Range("U3").Value = 2.25
$ws.Range("V3").Value = 1.57
$ws.Range("W3").Value = 5
$ws.Range("X3").Value = 6
$ws.Range("Y3").Value = 8.5
$ws.Range("Z3").Value = 10.5
$ws.Range("AA3").Value = 15
$ws.Range("AB3").Value = 37
$ws.Range("AD3").Value = 7
$ws.Range("AE3").Value = 22
$ws.Range("AF3").Value = 150
$ws.Range("AJ3").Value = 22
$ws.Range("AK3").Value = 175
$ws.Range("AL3").Value = 100
$ws.Range("AM3").Value = 90
$ws.Range("G4").Value = 2.37
$ws.Range("I4").Value = 3.55
$ws.Range("K4").Value = 1.87
$ws.Range("L4").Value = 4.15
$ws.Range("S4").Value = 1.53
$ws.Range("T4").Value = 2.32
$ws.Range("W4").Value = 6.1
$ws.Range("Z4").Value = 26
$ws.Range("AH4").Value = 7.9
$ws.Range("AI4").Value = 18
$ws.Range("AJ4").Value = 12
$ws.Range("AK4").Value = 55
$ws.Range("G5").Value = 1.7
$ws.Range("H5").Value = 3.4
$ws.Range("I5").Value = 5.1
$ws.Range("J5").Value = 2.2
$ws.Range("K5").Value = 2.15
$ws.Range("M5").Value = 1.34
$ws.Range("N5").Value = 3
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 1.72
$ws.Range("Q5").Value = 3.35
$ws.Range("R5").Value = 1.28
$ws.Range("S5").Value = 1.4
$ws.Range("T5").Value = 2.75
$ws.Range("U5").Value = 1.93
$ws.Range("V5").Value = 1.78
$ws.Range("W5").Value = 6.1
$ws.Range("AA5").Value = 14
$ws.Range("AB5").Value = 29
$ws.Range("AC5").Value = 6.7
$ws.Range("AE5").Value = 16.5
$ws.Range("AF5").Value = 90
$ws.Range("AG5").Value = 800
$ws.Range("AH5").Value = 12.5
$ws.Range("AI5").Value = 30
$ws.Range("AK5").Value = 100
$ws.Range("AM5").Value = 60
$ws.Range("AO5").Value = 6.7
$ws.Range("G8").Value = 2.15
$ws.Range("I8").Value = 3
$ws.Range("J8").Value = 2.75
$ws.Range("K8").Value = 2.4
$ws.Range("L8").Value = 3.25
$ws.Range("W8").Value = 13
$ws.Range("X8").Value = 15
$ws.Range("AI8").Value = 19
$ws.Range("AN8").Value = 1.02
$ws.Range("AO8").Value = 19
$ws.Range("G10").Value = 1.6
$ws.Range("H10").Value = 3.55
$ws.Range("I10").Value = 5.4
$ws.Range("J10").Value = 2.2
$ws.Range("K10").Value = 2.1
$ws.Range("L10").Value = 5.2
$ws.Range("N10").Value = 3.2
$ws.Range("Q10").Value = 2.77
$ws.Range("X10").Value = 7.5
$ws.Range("Z10").Value = 12
$ws.Range("AA10").Value = 13
$ws.Range("AD10").Value = 7.1
$ws.Range("AH10").Value = 15.5
$ws.Range("AJ10").Value = 16.5
$ws.Range("AL10").Value = 50
$ws.Range("G12").Value = 9
$ws.Range("H12").Value = 4.85
$ws.Range("I12").Value = 1.26
$ws.Range("J12").Value = 7.5
$ws.Range("K12").Value = 2.55
$ws.Range("L12").Value = 1.65
$ws.Range("P12").Value = 2.44
$ws.Range("Q12").Value = 2.12
$ws.Range("U12").Value = 1.9
$ws.Range("V12").Value = 1.81
$ws.Range("Y12").Value = 23
$ws.Range("AA12").Value = 80
$ws.Range("AB12").Value = 60
$ws.Range("AD12").Value = 8.75
$ws.Range("AE12").Value = 17
$ws.Range("AH12").Value = 7.1
$ws.Range("AI12").Value = 5.9
$ws.Range("G13").Value = 1.85
$ws.Range("H13").Value = 3.7
$ws.Range("I13").Value = 3.55
$ws.Range("J13").Value = 2.4
$ws.Range("L13").Value = 3.9
$ws.Range("N13").Value = 3.4
$ws.Range("Q13").Value = 2.6
$ws.Range("R13").Value = 1.38
$ws.Range("X13").Value = 9.5
$ws.Range("Z13").Value = 15.5
$ws.Range("AA13").Value = 14
$ws.Range("AC13").Value = 12.5
$ws.Range("AD13").Value = 7.3
$ws.Range("AF13").Value = 60
$ws.Range("AH13").Value = 12
$ws.Range("AI13").Value = 20
$ws.Range("AJ13").Value = 12.5
$ws.Range("AK13").Value = 50
$ws.Range("AL13").Value = 30
$ws.Range("AM13").Value = 35
$ws.Range("G14").Value = 1.67
$ws.Range("O14").Value = 2.1
$ws.Range("P14").Value = 1.7
$ws.Range("Q14").Value = 3.75
$ws.Range("R14").Value = 1.25
$ws.Range("G15").Value = 2.87
$ws.Range("I15").Value = 2.52
$ws.Range("J15").Value = 3.5
$ws.Range("K15").Value = 1.95
$ws.Range("L15").Value = 3.2
$ws.Range("M15").Value = 1.44
$ws.Range("N15").Value = 2.6
$ws.Range("S15").Value = 1.5
$ws.Range("T15").Value = 2.42
$ws.Range("U15").Value = 1.93
$ws.Range("V15").Value = 1.78
$ws.Range("W15").Value = 7.5
$ws.Range("X15").Value = 13.5
$ws.Range("Y15").Value = 10.75
$ws.Range("Z15").Value = 35
$ws.Range("AA15").Value = 28
$ws.Range("AH15").Value = 6.7
$ws.Range("AI15").Value = 11.25
$ws.Range("AJ15").Value = 9.75
$ws.Range("AK15").Value = 27
$ws.Range("AL15").Value = 24
$ws.Range("G17").Value = 3.2
$ws.Range("H17").Value = 2.65
$ws.Range("I17").Value = 2.65
$ws.Range("L17").Value = 3.6
$ws.Range("M17").Value = 1.63
$ws.Range("N17").Value = 2.1
$ws.Range("R17").Value = 1.08
$ws.Range("AI17").Value = 11
$ws.Range("AN17").Value = 1.13
$ws.Range("AO17").Value = 5
$ws.Range("G18").Value = 1.92
$ws.Range("H18").Value = 3.5
$ws.Range("M18").Value = 1.25
$ws.Range("O18").Value = 1.93
$ws.Range("P18").Value = 1.88
$ws.Range("Q18").Value = 3.25
$ws.Range("AA18").Value = 15
$ws.Range("AC18").Value = 11
$ws.Range("AD18").Value = 7
$ws.Range("AN18").Value = 1.03
$ws.Range("O19").Value = 2.08
$ws.Range("P19").Value = 1.73
$ws.Range("AH19").Value = 11
$ws.Range("M20").Value = 1.36
$ws.Range("N20").Value = 3
$ws.Range("O21").Value = 1.93
$ws.Range("P21").Value = 1.88
$ws.Range("Q21").Value = 3.25
$ws.Range("R21").Value = 1.33
$ws.Range("AN21").Value = 1.05
$ws.Range("AO21").Value = 11
